$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Federico Fasanelli"
$ws.Range("B17").Value = "Stefano Tita | Clitoriders"
$ws.Range("C17").Value = "Nadir Chtioui | MAI UNA GIOIA"
$ws.Range("D17").Value = "Federico  Manica | iMontagna"
$ws.Range("E17").Value = "Mattia Festi | Shark Attack"
$ws.Range("F17").Value = "Alessandro  Tengattini | Herta Vernello"
